$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append new data row for DGS's 2021/09/20 report
$row = 83
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2021/09/20"
$cellA.NumberFormat = "yyyy/mm/dd"
$ws.Cells.Item($row, 2).Value = 149.1
$ws.Cells.Item($row, 3).Value = 152.4
$ws.Cells.Item($row, 4).Value = 0.82
$ws.Cells.Item($row, 5).Value = 0.81

# Mirror the author's UI state: active cell moves to the next empty row
[void]$ws.Range("A84").Select()
